$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need to be forced to Text
# format first, otherwise Excel auto-converts them to numeric values and we
# lose exact formatting (e.g. trailing zeros like "1.00" -> 1).
$textForceCells = @(
    "D5",
    "D6",
    "D7",
    "D12",
    "D14",
    "D16",
    "D18",
    "D19",
    "D21",
    "D22",
    "D23",
    "D24",
    "D25",
    "D26",
    "D27",
    "D28",
    "D29",
    "D30",
    "D32",
    "D33",
    "D35",
    "D36",
    "D37",
    "D38",
    "D39",
    "D43",
    "D45",
    "D47",
    "D48",
    "D49",
    "D50",
    "D51"
)
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated cell values from the crypto data refresh
$ws.Range("D2").Value2 = '63.672.76'
$ws.Range("E2").Value2 = '  -3.28%  '
$ws.Range("D3").Value2 = '3.323.93'
$ws.Range("E3").Value2 = '  -4.80%  '
$ws.Range("E4").Value2 = '  +0.23%  '
$ws.Range("D5").Value2 = '548.04'
$ws.Range("E5").Value2 = '  -1.81%  '
$ws.Range("D6").Value2 = '171.85'
$ws.Range("E6").Value2 = '  -5.20%  '
$ws.Range("D7").Value2 = '0.612'
$ws.Range("E7").Value2 = '  -4.18%  '
$ws.Range("E8").Value2 = '  -0.05%  '
$ws.Range("D9").Value2 = '3.314.07'
$ws.Range("E9").Value2 = '  -4.87%  '
$ws.Range("E10").Value2 = '  -4.61%  '
$ws.Range("E11").Value2 = '  -2.17%  '
$ws.Range("D12").Value2 = '53.36'
$ws.Range("E12").Value2 = '  -2.37%  '
$ws.Range("E13").Value2 = '  -3.54%  '
$ws.Range("B14").Value2 = 'Polkadot'
$ws.Range("C14").Value2 = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").Value2 = '8.86'
$ws.Range("E14").Value2 = '  -4.99%  '
$ws.Range("B15").Value2 = 'WrappedliquidstakedEther2.0'
$ws.Range("C15").Value2 = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D15").Value2 = '3.871.55'
$ws.Range("E15").Value2 = '  -3.55%  '
$ws.Range("D16").Value2 = '18.25'
$ws.Range("E16").Value2 = '  -2.86%  '
$ws.Range("D17").Value2 = '3.337.28'
$ws.Range("E17").Value2 = '  -3.66%  '
$ws.Range("D18").Value2 = '0.117'
$ws.Range("E18").Value2 = '  -3.68%  '
$ws.Range("D19").Value2 = '11.64'
$ws.Range("E19").Value2 = '  -3.80%  '
$ws.Range("D20").Value2 = '63.550.73'
$ws.Range("E20").Value2 = '  -3.28%  '
$ws.Range("D21").Value2 = '0.974'
$ws.Range("E21").Value2 = '  -2.20%  '
$ws.Range("D22").Value2 = '410.03'
$ws.Range("E22").Value2 = '  -2.46%  '
$ws.Range("B23").Value2 = 'Toncoin'
$ws.Range("C23").Value2 = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D23").Value2 = '4.39'
$ws.Range("E23").Value2 = '  +6.30%  '
$ws.Range("B24").Value2 = 'PancakeSwap'
$ws.Range("C24").Value2 = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D24").Value2 = '4.03'
$ws.Range("E24").Value2 = '  -1.21%  '
$ws.Range("D25").Value2 = '13.76'
$ws.Range("E25").Value2 = '  +7.59%  '
$ws.Range("D26").Value2 = '82.82'
$ws.Range("E26").Value2 = '  -4.54%  '
$ws.Range("D27").Value2 = '10.50'
$ws.Range("E27").Value2 = '  -3.97%  '
$ws.Range("D28").Value2 = '2.72'
$ws.Range("E28").Value2 = '  -6.01%  '
$ws.Range("D29").Value2 = '8.59'
$ws.Range("E29").Value2 = '  -5.97%  '
$ws.Range("D30").Value2 = '28.98'
$ws.Range("E30").Value2 = '  -5.20%  '
$ws.Range("E31").Value2 = '  -4.12%  '
$ws.Range("D32").Value2 = '578.02'
$ws.Range("E32").Value2 = '  -5.79%  '
$ws.Range("D33").Value2 = '11.32'
$ws.Range("E33").Value2 = '  -4.33%  '
$ws.Range("E34").Value2 = '  -4.30%  '
$ws.Range("D35").Value2 = '57.81'
$ws.Range("E35").Value2 = '  -2.85%  '
$ws.Range("B36").Value2 = 'Kaspa'
$ws.Range("C36").Value2 = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D36").Value2 = '0.147'
$ws.Range("E36").Value2 = '  +0.70%  '
$ws.Range("B37").Value2 = 'Dai'
$ws.Range("C37").Value2 = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D37").Value2 = '1.00'
$ws.Range("E37").Value2 = '  -0.08%  '
$ws.Range("D38").Value2 = '34.91'
$ws.Range("E38").Value2 = '  -7.81%  '
$ws.Range("D39").Value2 = '3.39'
$ws.Range("E39").Value2 = '  +1.95%  '
$ws.Range("D40").Value2 = '0.0₃0734'
$ws.Range("E40").Value2 = '  -8.25%  '
$ws.Range("E41").Value2 = '  -5.32%  '
$ws.Range("D42").Value2 = '3.117.54'
$ws.Range("E42").Value2 = '  -6.45%  '
$ws.Range("D43").Value2 = '0.998'
$ws.Range("E43").Value2 = '  +0.40%  '
$ws.Range("E44").Value2 = '  -2.51%  '
$ws.Range("D45").Value2 = '3.22'
$ws.Range("E45").Value2 = '  -1.07%  '
$ws.Range("E46").Value2 = '  -4.70%  '
$ws.Range("D47").Value2 = '2.41'
$ws.Range("E47").Value2 = '  -6.29%  '
$ws.Range("D48").Value2 = '2.60'
$ws.Range("E48").Value2 = '  -4.26%  '
$ws.Range("D49").Value2 = '0.128'
$ws.Range("E49").Value2 = '  -4.15%  '
$ws.Range("D50").Value2 = '132.79'
$ws.Range("E50").Value2 = '  -3.94%  '
$ws.Range("D51").Value2 = '8.00'
$ws.Range("E51").Value2 = '  -6.08%  '

